$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new test row for Pain_Index macro (write in C, B, A order so the
# shared-string table is populated in the same order Excel produced it)
$ws.Range("C86").Value = "Pain_Index_test"
$ws.Range("B86").Value = "Test pain index"
$ws.Range("A86").Value = "Pain index"

# Update the active selection to match the new last cell used
$ws.Range("A86").Select()
